# Add new COVID-19 data rows for Pakistan (2021-03-11 .. 2021-03-15),
# appended after the existing last row (381) on the "Covid-19" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19")

# Columns: A=Country, B=Report_Date, C=Confirmed, D=Deaths, E=Recovered,
#          F=Active, G=New Confirmed (=C-C_prev), H=New Deaths (=D-D_prev),
#          I=New Recovered (=E-E_prev)
$newRows = @(
    @{ Row = 382; Date = 44266; Confirmed = 600198; Deaths = 13430; Recovered = 568065; Active = 18703 },
    @{ Row = 383; Date = 44267; Confirmed = 602536; Deaths = 13476; Recovered = 569296; Active = 19764 },
    @{ Row = 384; Date = 44268; Confirmed = 605200; Deaths = 13508; Recovered = 570571; Active = 21121 },
    @{ Row = 385; Date = 44269; Confirmed = 607453; Deaths = 13537; Recovered = 571878; Active = 22038 },
    @{ Row = 386; Date = 44270; Confirmed = 609964; Deaths = 13595; Recovered = 573014; Active = 23355 }
)

# Seed formatting for the new rows by copying the previous (last existing)
# row's formats down, same as Excel does when a user fills a new row under
# an existing, formatted table.
$ws.Range("A381:I381").Copy() | Out-Null
$ws.Range("A382:I386").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = "Pakistan"
    $ws.Cells.Item($row, 2).Value = $r.Date

    $ws.Cells.Item($row, 3).Value = $r.Confirmed
    $ws.Cells.Item($row, 4).Value = $r.Deaths
    $ws.Cells.Item($row, 5).Value = $r.Recovered
    $ws.Cells.Item($row, 6).Value = $r.Active

    $ws.Cells.Item($row, 7).Formula = "=C$row-C$($row-1)"
    $ws.Cells.Item($row, 8).Formula = "=D$row-D$($row-1)"
    $ws.Cells.Item($row, 9).Formula = "=E$row-E$($row-1)"
}

$ws.Range("I389").Select()
$excel.ActiveWindow.ScrollRow = 375

$wb.Application.Calculate()
